$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the PID values (column C, rows 2-26) while keeping formatting/styles.
$ws.Range("C2:C26").ClearContents()

# Match the new active selection recorded in the saved file.
$ws.Range("C2:C26").Select()
